$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" — the handoff run picked up the
# 6353d0f2-...md file and moved it from "In Translation" to
# "Ready for handoff" (stamping new handoff timestamps). The report rows are
# re-sorted by file name, so c53cee36-...md (alphabetically first) now sits
# on row 2 and 6353d0f2-...md drops to row 3 on every sheet.
# ---------------------------------------------------------------------------

$urlMdC53   = "https://github.com/OpenLocalizationTest/oltest/blob/9d8422620be65b5c029c87d87131d84e5456340d/e2e/c53cee36-a985-48c3-8281-f654ea4f7aba.md"
$urlMd6353  = "https://github.com/OpenLocalizationTest/oltest/blob/9d8422620be65b5c029c87d87131d84e5456340d/e2e/6353d0f2-2090-4abf-8edc-4a1f24ea9957.md"
$urlZhC53   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bca44305e8dbb270144ee4cdd3f2aef1704801e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c53cee36-a985-48c3-8281-f654ea4f7aba.4037a86397cd9ead38c8ea6a3a3acbf621167eec.zh-cn.xlf"
$urlZh6353  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bca44305e8dbb270144ee4cdd3f2aef1704801e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6353d0f2-2090-4abf-8edc-4a1f24ea9957.9182ccc28a4623335c6d0135dab6dfdffcd7a536.zh-cn.xlf"
$urlDeC53   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82ef975e2c0a3e5d9dca84309edacc6575a85efc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c53cee36-a985-48c3-8281-f654ea4f7aba.4037a86397cd9ead38c8ea6a3a3acbf621167eec.de-de.xlf"
$urlDe6353  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82ef975e2c0a3e5d9dca84309edacc6575a85efc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6353d0f2-2090-4abf-8edc-4a1f24ea9957.9182ccc28a4623335c6d0135dab6dfdffcd7a536.de-de.xlf"

$mdC53  = "c53cee36-a985-48c3-8281-f654ea4f7aba.md"
$md6353 = "6353d0f2-2090-4abf-8edc-4a1f24ea9957.md"
$zhC53  = "c53cee36-a985-48c3-8281-f654ea4f7aba.4037a86397cd9ead38c8ea6a3a3acbf621167eec.zh-cn.xlf"
$zh6353 = "6353d0f2-2090-4abf-8edc-4a1f24ea9957.9182ccc28a4623335c6d0135dab6dfdffcd7a536.zh-cn.xlf"
$deC53  = "c53cee36-a985-48c3-8281-f654ea4f7aba.4037a86397cd9ead38c8ea6a3a3acbf621167eec.de-de.xlf"
$de6353 = "6353d0f2-2090-4abf-8edc-4a1f24ea9957.9182ccc28a4623335c6d0135dab6dfdffcd7a536.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
$wsOverview.Range("D2").Value = "2016-03-23 14:16:57"

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-23 14:17:58"

$wsOverview.Hyperlinks.Item(1).Delete()
$wsOverview.Hyperlinks.Item(1).Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $urlMdC53, "", "", $mdC53)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urlMd6353, "", "", $md6353)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("E2").Value = "2016-03-23 14:16:53"

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-23 14:17:54"

$wsZhCn.Hyperlinks.Item(1).Delete()
$wsZhCn.Hyperlinks.Item(1).Delete()
$wsZhCn.Hyperlinks.Item(1).Delete()
$wsZhCn.Hyperlinks.Item(1).Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlMdC53, "", "", $mdC53)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $urlZhC53, "", "", $zhC53)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlMd6353, "", "", $md6353)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $urlZh6353, "", "", $zh6353)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("E2").Value = "2016-03-23 14:16:57"

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-23 14:17:58"

$wsDeDe.Hyperlinks.Item(1).Delete()
$wsDeDe.Hyperlinks.Item(1).Delete()
$wsDeDe.Hyperlinks.Item(1).Delete()
$wsDeDe.Hyperlinks.Item(1).Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlMdC53, "", "", $mdC53)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $urlDeC53, "", "", $deC53)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlMd6353, "", "", $md6353)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $urlDe6353, "", "", $de6353)
